$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force every target cell to Text format first so ambiguous numeric-looking
# strings (e.g. "0.999", "1.00", "544.43") are preserved verbatim as text,
# matching the original inlineStr cells instead of being parsed into floats.
$cells = @("D2", "D3", "E3", "D4", "D5", "E5", "D6", "E6", "E7", "E8", "D9", "E9", "E10", "D11", "E11", "E12", "D13", "E13", "D14", "D15", "D16", "E16", "E17", "D18", "E18", "E19", "D20", "E20", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "E25", "D26", "E26", "D27", "E27", "B28", "C28", "D28", "E28", "B29", "C29", "D29", "E29", "D30", "E30", "E31", "E32", "D33", "E33", "E34", "D35", "E35", "E36", "E37", "D38", "E38", "D39", "E39", "E40", "D41", "E41", "E42", "E43", "E44", "D45", "E45", "E46", "E47", "D48", "E48", "D49", "E49", "E50", "D51", "E51")
foreach ($cell in $cells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "59.340.82"
$ws.Range("D3").Value = "2.608.44"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").Value = "544.43"
$ws.Range("E5").Value = "  +4.39%  "
$ws.Range("D6").Value = "141.23"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "6.46"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").Value = "0.335"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("D13").Value = "3.067.05"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").Value = "59.257.09"
$ws.Range("D15").Value = "20.58"
$ws.Range("D16").Value = "2.616.65"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "343.93"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").Value = "10.13"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "67.55"
$ws.Range("E23").Value = "  +2.18%  "
$ws.Range("D24").Value = "0.409"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("E25").Value = "  -1.02%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("D27").Value = "7.22"
$ws.Range("E27").Value = "  +1.31%  "
$ws.Range("B28").Value = "USDe"
$ws.Range("C28").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0739"
$ws.Range("E29").Value = "  +1.62%  "
$ws.Range("D30").Value = "1.72"
$ws.Range("E30").Value = "  +9.59%  "
$ws.Range("E31").Value = "  -2.50%  "
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("D33").Value = "149.52"
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("E34").Value = "  -0.62%  "
$ws.Range("D35").Value = "37.16"
$ws.Range("E35").Value = "  +2.18%  "
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("D38").Value = "0.836"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").Value = "0.815"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("E40").Value = "  +1.02%  "
$ws.Range("D41").Value = "277.10"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("E43").Value = "  +1.25%  "
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").Value = "0.0956"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("E47").Value = "  +1.34%  "
$ws.Range("D48").Value = "1.943.98"
$ws.Range("E48").Value = "  -2.17%  "
$ws.Range("D49").Value = "18.48"
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("E50").Value = "  -2.51%  "
$ws.Range("D51").Value = "111.15"
$ws.Range("E51").Value = "  -2.13%  "

# Restore default (Normal) style so no stray number-format/style is left
# behind on cells that should remain visually identical to before.
foreach ($cell in $cells) {
    $ws.Range($cell).Style = "Normal"
}
